$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on columns that would otherwise be auto-converted
# by Excel's input parser (dates in A, "$"-prefixed amounts in B), so the
# values stay literal strings matching the source.
$ws.Range("A2:A12").NumberFormat = "@"
$ws.Range("B2:B12").NumberFormat = "@"

# Update existing row 2 values (B2, C2 change; E2 gets a new value)
$ws.Range("B2").Value = "-`$3"
$ws.Range("C2").Value = "store 1"
$ws.Range("E2").Value = "-"

# Append new rows 3-12 as text values (matches inlineStr semantics of source)
$data = @(
    @("9/17/2022",  "-`$53.89", "Target",           "Ximena Leyva", "Drinks for party"),
    @("9/27/2022",  "-`$43.78", "another store",    "Ximena Leyva", "-"),
    @("10/8/2022",  "-`$73",    "a place",          "Ximena Leyva", "-"),
    @("10/27/2022", "-`$142",   "shoe store",       "Ximena Leyva", "shoes"),
    @("11/8/2022",  "-`$34",    "make up store",    "Ximena Leyva", "-"),
    @("10/5/2022",  "+`$430",   "funding 2",        "Ximena Leyva", "-"),
    @("11/10/2022", "+`$534",   "source 4",         "Ximena Leyva", "-"),
    @("9/9/2022",   "+`$24",    "donation",         "Ximena Leyva", "-"),
    @("11/1/2022",  "+`$1000",  "CPA",              "Ximena Leyva", "-"),
    @("12/18/2022", "+`$2500",  "massive donation", "Ximena Leyva", "-")
)

$row = 3
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
